$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data change: only the Products suite should run now ---
# Runmode column: Customer -> N, MyAccount -> N, Products stays Y
$ws.Range("C2").Value2 = "N"
$ws.Range("C3").Value2 = "N"
$ws.Range("C4").Value2 = "Y"

# --- Formatting: give the table a clean bordered look ---
$full = $ws.Range("A1:C4")

# Thin grid lines on the inside of the table
$full.Borders.Item(11).LineStyle = 1   # xlInsideVertical
$full.Borders.Item(11).Weight = 2      # xlThin
$full.Borders.Item(12).LineStyle = 1   # xlInsideHorizontal
$full.Borders.Item(12).Weight = 2      # xlThin

# Medium (thicker) box border around the whole table
$full.BorderAround(1, -4138)           # xlContinuous, xlMedium

# Medium separator between the header row and the data rows
$sep = $ws.Range("A1:C2")
$sep.Borders.Item(12).Weight = -4138   # xlInsideHorizontal, xlMedium

# Center the header row text
$headerRow = $ws.Range("A1:C1")
$headerRow.HorizontalAlignment = -4108 # xlCenter

# Row heights: header row grows slightly, last row grows slightly (thicker border)
$ws.Rows.Item(1).RowHeight = 16
$ws.Rows.Item(4).RowHeight = 106

# --- Selection moves to E3 ---
$ws.Range("E3").Select()
